# Refresh the cryptocurrency "Price" (column D) and "Volume(1h)" change
# (column E) figures on the active worksheet with the latest pulled values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.043.23'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = '''1.800.48'
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''309.92'
$ws.Range("E5").Value = '  -1.49%  '

$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '''0.5076'
$ws.Range("E7").Value = '  -4.15%  '

$ws.Range("D8").Value = '''0.3829'
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.07719'
$ws.Range("E9").Value = '  -3.49%  '

$ws.Range("D10").Value = '''1.093'
$ws.Range("E10").Value = '  -0.62%  '

$ws.Range("E11").Value = '  -1.61%  '

$ws.Range("D12").Value = '''6.357'
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").Value = '''1.003'
$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").Value = '''20.33'
$ws.Range("E14").Value = '  -1.34%  '

$ws.Range("D15").Value = '''1.801.88'
$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("D16").Value = '''7.268'
$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").Value = '''0.00001070'
$ws.Range("E18").Value = '  -2.39%  '

$ws.Range("D19").Value = '''0.06567'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").Value = '''17.23'
$ws.Range("E21").Value = '  -0.82%  '

$ws.Range("D22").Value = '''5.978'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").Value = '''28.051.74'

$ws.Range("D24").Value = '''11.03'
$ws.Range("E24").Value = '  -1.63%  '

$ws.Range("D25").Value = '''2.223'
$ws.Range("E25").Value = '  -0.50%  '

$ws.Range("D26").Value = '''159.33'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").Value = '''2.416'
$ws.Range("E27").Value = '  +1.39%  '

$ws.Range("D28").Value = '''2.011.48'
$ws.Range("E28").Value = '  +0.16%  '

$ws.Range("D29").Value = '''20.21'

$ws.Range("D30").Value = '''127.10'
$ws.Range("E30").Value = '  +3.12%  '

$ws.Range("E31").Value = '  +0.08%  '

$ws.Range("E32").Value = '  -1.50%  '

$ws.Range("D33").Value = '''3.647'
$ws.Range("E33").Value = '  -0.49%  '

$ws.Range("D34").Value = '''5.535'
$ws.Range("E34").Value = '  -0.37%  '

$ws.Range("D35").Value = '''0.06942'
$ws.Range("E35").Value = '  -4.64%  '

$ws.Range("D36").Value = '''9.099'
$ws.Range("E36").Value = '  +1.91%  '

$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("D38").Value = '''0.2165'
$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("D39").Value = '''4.991'
$ws.Range("E39").Value = '  -1.73%  '

$ws.Range("E40").Value = '  -6.98%  '

$ws.Range("D41").Value = '''0.6097'
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("E42").Value = '  +0.14%  '

$ws.Range("D43").Value = '''1.149'
$ws.Range("E43").Value = '  -1.63%  '

$ws.Range("D44").Value = '''13.24'
$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("E45").Value = '  -5.87%  '

$ws.Range("E46").Value = '  -2.06%  '

$ws.Range("D47").Value = '''3.703'
$ws.Range("E47").Value = '  -1.67%  '

$ws.Range("D48").Value = '''125.73'
$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("D49").Value = '''1.928'
$ws.Range("E49").Value = '  -0.08%  '

$ws.Range("D50").Value = '''1.182'
$ws.Range("E50").Value = '  -2.44%  '

$ws.Range("D51").Value = '''0.06720'
$ws.Range("E51").Value = '  -1.67%  '
